$d = $word.ActiveDocument

$d.Content.Find.Execute("day to day tasks", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "day-to-day tasks", 2)

$d.Content.Find.Execute("AI, and level while leading", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "AI, and levels while leading", 2)

$d.Content.Find.Execute("a  first person game", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "a  first-person game", 2)

$d.Content.Find.Execute("a direct focus in game design.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "a direct focus on game design.", 2)
